$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Financial Aid" table (G1:H7) renamed to "Scholarships" table ---
$ws.Range("G2").Value = "Scholarships"
$ws.Range("G3").Value = "sholarship name"
$ws.Range("G5").Value = "amount"
$ws.Range("H6").Value = "scholarship"

# The rows that previously had the highlighted "amount"/"tuition" fill
# (G4:H4, H5) lose their special fill now that the table no longer has a
# distinct "amount" row further down - restore them to the default style.
$ws.Range("G4:H4").ClearFormats()
$ws.Range("H5").ClearFormats()

# Update the active selection to reflect where the author ended up working
$ws.Range("I12").Select()
